# Weekly update: insert the newest "Femacal de La Calera - Breva" price
# record at the top of the data block (row 9), pushing the existing rows
# down by one (old row 9 -> 10, old row 10 -> 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 9 and 10 down to 10 and 11.
$ws.Rows("9:9").Insert()

# Populate the new row 9 with this week's record.
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "Femacal de La Calera"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44907
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101006
$ws.Range("J9").Value = "Breva"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 25000
$ws.Range("Q9").Value = "$/bandeja 5 kilos"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 5000
$ws.Range("T9").Value = 5
